$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before AT. This shifts the cell VALUES of AT:BB one
# column to the right (AT->AU, ..., BB->BC) but, unlike real Excel, it does
# NOT move the cell comments along with them - comments stay attached to
# their original column letters. So after the value-shifting Insert(), walk
# the old AT:BB range from right to left and re-home each comment onto the
# column immediately to its right to match.
$ws.Columns("AT:AT").Insert()

$srcCols = @("AT","AU","AV","AW","AX","AY","AZ","BA","BB")
$dstCols = @("AU","AV","AW","AX","AY","AZ","BA","BB","BC")
for ($i = $srcCols.Length - 1; $i -ge 0; $i--) {
    $src = $srcCols[$i] + "15"
    $dst = $dstCols[$i] + "15"
    $commentText = $ws.Range($src).Comment.Text()
    $ws.Range($src).Comment.Delete()
    $ws.Range($dst).AddComment($commentText)
}

# Populate the new AT15 header cell with the specimen_voucher text and
# attach its descriptive comment.
$ws.Range("AT15").Value = "specimen_voucher"
$ws.Range("AT15").AddComment("Identifier for the physical specimen. Use format: ""[<institution-code>:[<collection-code>:]]<specimen_id>"", eg, ""UAM:Mamm:52179"". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a 'structured voucher'. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier")
